# ALS application form: closing date changed + "_GoBack" bookmark relocated
# (the bookmark move / stray empty-paragraph cleanup are side effects of the
# author re-opening and re-saving the doc in Word; the real content change is
# the closing-date sentence).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the very start of the
#    document (it sat around an empty range on the first paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. First "key criteria" table: drop the trailing blank paragraph
#    (cell had 5 empty justified paragraphs, now has 4).
# ------------------------------------------------------------------
$table1 = $d.Tables.Item(3)
$cell1 = $table1.Cell(1, 1)
$cell1.Range.Paragraphs.Item(5).Range.Delete()

# ------------------------------------------------------------------
# 3. Second "key criteria" table: drop the 2nd blank paragraph and
#    re-add the "_GoBack" bookmark around the (now first remaining)
#    blank paragraph.
# ------------------------------------------------------------------
$table2 = $d.Tables.Item(4)
$cell2 = $table2.Cell(1, 1)
$d.Bookmarks.Add("_GoBack", $cell2.Range.Paragraphs.Item(1).Range) | Out-Null
$cell2.Range.Paragraphs.Item(2).Range.Delete()

# ------------------------------------------------------------------
# 4. Update the closing date sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Friday 29 January. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Monday 1st February (midnight). ", 2)
